$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Phase 1: seed shared-string insertion order to match target sharedStrings.xml ---
# (one "first touch" cell per new unique string, in the exact order the strings
#  must appear in xl/sharedStrings.xml)
$ws.Range("K130").Value = "abcd"
$ws.Range("S126").Value = "a"
$ws.Range("S127").Value = "b"
$ws.Range("T127").Value = "c"
$ws.Range("T126").Value = "d"
$ws.Range("L130").Value = "dcba"
$ws.Range("M130").Value = "lls"
$ws.Range("N130").Value = "s"
$ws.Range("O130").Value = "sssll"
$ws.Range("U126").Value = "l"
$ws.Range("U133").Value = "Trie"
$ws.Range("K156").Value = "ab"
$ws.Range("K157").Value = "ac"
$ws.Range("K158").Value = "aa"

# --- Phase 2: numeric cells (row 125 header: 0..4) ---
$ws.Range("S125").Value = 0
$ws.Range("T125").Value = 1
$ws.Range("U125").Value = 2
$ws.Range("V125").Value = 3
$ws.Range("W125").Value = 4

# --- Phase 3: remaining string cells for the new Trie / Palindrome-Pairs tables ---
$ws.Range("V126").Value = "s"
$ws.Range("W126").Value = "s"
$ws.Range("U127").Value = "l"
$ws.Range("W127").Value = "s"
$ws.Range("S128").Value = "c"
$ws.Range("T128").Value = "b"
$ws.Range("U128").Value = "s"
$ws.Range("W128").Value = "s"
$ws.Range("S129").Value = "d"
$ws.Range("T129").Value = "a"
$ws.Range("W129").Value = "l"
$ws.Range("W130").Value = "l"

$ws.Range("K134").Value = "lls"
$ws.Range("L134").Value = "abcd"
$ws.Range("T134").Value = "a"
$ws.Range("U134").Value = "d"
$ws.Range("V134").Value = "l"
$ws.Range("W134").Value = "s"

$ws.Range("K135").Value = "lls"
$ws.Range("L135").Value = "dcba"
$ws.Range("T135").Value = "b"
$ws.Range("U135").Value = "c"
$ws.Range("V135").Value = "l"
$ws.Range("W135").Value = "s"

$ws.Range("K136").Value = "lls"
$ws.Range("L136").Value = "lls"
$ws.Range("T136").Value = "c"
$ws.Range("U136").Value = "b"
$ws.Range("V136").Value = "s"
$ws.Range("W136").Value = "s"

$ws.Range("K137").Value = "lls"
$ws.Range("L137").Value = "s"
$ws.Range("T137").Value = "d"
$ws.Range("U137").Value = "a"
$ws.Range("W137").Value = "l"

$ws.Range("K138").Value = "lls"
$ws.Range("L138").Value = "sssll"
$ws.Range("W138").Value = "l"

$ws.Range("K140").Value = "s"
$ws.Range("L140").Value = "abcd"
$ws.Range("O140").Value = "abcd"
$ws.Range("Q140").Value = "abcd"
$ws.Range("R140").Value = "dcba"

$ws.Range("K141").Value = "s"
$ws.Range("L141").Value = "dcba"
$ws.Range("O141").Value = "dcba"
$ws.Range("Q141").Value = "dcba"
$ws.Range("R141").Value = "abcd"

$ws.Range("K142").Value = "s"
$ws.Range("L142").Value = "lls"
$ws.Range("O142").Value = "s"
$ws.Range("Q142").Value = "s"
$ws.Range("R142").Value = "lls"

$ws.Range("K143").Value = "s"
$ws.Range("L143").Value = "s"
$ws.Range("O143").Value = "lls"
$ws.Range("Q143").Value = "lls"
$ws.Range("R143").Value = "sssll"

$ws.Range("K144").Value = "s"
$ws.Range("L144").Value = "sssll"

$ws.Range("S152").Value = "Trie"

$ws.Range("K153").Value = "a"
$ws.Range("R153").Value = "a"
$ws.Range("T153").Value = "b"
$ws.Range("U153").Value = "c"

$ws.Range("K154").Value = "b"

$ws.Range("K155").Value = "c"
$ws.Range("Q155").Value = "b"
$ws.Range("R155").Value = "c"
$ws.Range("S155").Value = "a"

# --- Phase 4: formulas (CONCAT of L & K columns) ---
$ws.Range("M134").Formula = "=CONCAT(L134,K134)"
$ws.Range("M135").Formula = "=CONCAT(L135,K135)"
$ws.Range("M137").Formula = "=CONCAT(L137,K137)"
$ws.Range("M138").Formula = "=CONCAT(L138,K138)"
$ws.Range("M139").Formula = "=CONCAT(L139,K139)"
$ws.Range("M140").Formula = "=CONCAT(L140,K140)"
$ws.Range("M141").Formula = "=CONCAT(L141,K141)"
$ws.Range("M142").Formula = "=CONCAT(L142,K142)"
$ws.Range("M144").Formula = "=CONCAT(L144,K144)"

# --- Phase 5: merge U133:V133 ("Trie" header) and center it ---
$ws.Range("U133:V133").Merge()
$ws.Range("U133:V133").HorizontalAlignment = -4108
$ws.Range("U133:V133").VerticalAlignment = -4108

# --- Phase 6: widen column M (bestfit for the longer CONCAT results) ---
$ws.Columns.Item(13).AutoFit()

# --- Phase 7: update view / selection to match target workbook state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 148
$win.ScrollColumn = 9
$ws.Range("V154").Select()
